# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Column G on Sheet1 holds "K" (number of strikes / breaks-style count used by
# the downstream std/mean calculation). The values below are the freshly
# recomputed s_vals for each row (r = 2..45), replacing the previous
# "Strike#"-derived numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 1
    8  = 0
    9  = 1
    10 = 1
    11 = 1
    12 = 1
    13 = 2
    14 = 1
    15 = 0
    16 = 2
    17 = 1
    18 = 1
    19 = 0
    20 = 2
    21 = 0
    22 = 1
    23 = 2
    24 = 0
    25 = 2
    26 = 0
    27 = 3
    28 = 2
    29 = 0
    30 = 4
    31 = 1
    32 = 0
    33 = 1
    34 = 3
    35 = 3
    36 = 3
    37 = 0
    38 = 3
    39 = 3
    40 = 1
    41 = 5
    42 = 4
    43 = 4
    44 = 3
    45 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
